# CTP review round 2 update for tobaccoproblem_reference.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: report_id ---
# Type column: integer -> number (description unchanged)
$ws.Range("C2").Value = "number"

# --- Row 3: submission_date ---
# Description: "Date of Original Receipt to FDA." -> new longer text about CTP receipt
$ws.Range("D3").Value = "Date report was received by CTP; this is the earliest date of report receipt, either to Safety Reporting Portal (SRP) or by other means"
$ws.Rows(3).RowHeight = 34

# --- Row 4: reported_product_problems ---
$ws.Rows(4).RowHeight = 17

# --- Row 5: reported_health_problems ---
$ws.Rows(5).RowHeight = 17

# --- Row 6: product_type ---
$ws.Rows(6).RowHeight = 17

# --- Row 7: number_tobacco_products ---
# Type column: integer -> number
$ws.Range("C7").Value = "number"

# --- Row 8: number_product_problems ---
# Type column: integer -> number
$ws.Range("C8").Value = "number"

# --- Row 9: number_health_problems ---
# Type column: integer -> number; Description: fix typo "ystem-calculated" -> "System-calculated"
$ws.Range("C9").Value = "number"
$ws.Range("D9").Value = "System-calculated number of Health Problems (i.e., MedDRA terms selected from a standardized list of symptoms, signs, diagnoses and outcomes) reported, displayed as a whole number, ≥0."

# --- Row 10: nonuser_affected ---
# Description updated with date-range clause
$ws.Range("D10").Value = "Displays text reflecting the response to this optional question (2017 - 12/14/2018) or required question (12/15/2018 onward) as “No information provided” if not answered, or Yes/No."
$ws.Rows(10).RowHeight = 51

# --- Sheet view / print setup ---
$null = $ws.Range("F10").Select()
$ws.PageSetup.Zoom = 65
